$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Dep or Non-Dep Without LOSOCV")
$ws2 = $wb.Worksheets.Item("Dep or Non-Dep With LOSOCV")

# 1. Update the two title cells (drop the word "Results" from the titles).
$ws1.Range("A1").Value = "Machine Learning - Depression/Non-Depression Classification (without Leave One Subject Out Cross Validation)"
$ws2.Range("A1").Value = "Machine Learning - Depression/Non-Depression Classification (with Leave One Subject Out Cross Validation)"

# 2. Fill in the new K-Nearest Neighbours results rows (row 6 and row 13) on both sheets.

# --- Sheet 1: "Dep or Non-Dep Without LOSOCV" ---
$ws1.Range("B6").Value = 0.89156626506024095
$ws1.Range("C6").Value = 0.96296296296296202
$ws1.Range("D6").Value = 0.76470588235294101
$ws1.Range("E6").Value = 0.85245901639344202
$ws1.Range("F6").Value = 0.87214885954381705
$ws1.Range("F6").NumberFormat = "0.000000"

# Row 12 (Random Forest) keeps its existing values, but B12 picks up the
# six-decimal number format in the source edit as well.
$ws1.Range("B12").NumberFormat = "0.000000"

$ws1.Range("B13").Value = 0.72289156626506001
$ws1.Range("C13").Value = 0.73913043478260798
$ws1.Range("D13").Value = 0.5
$ws1.Range("E13").Value = 0.59649122807017496
$ws1.Range("F13").Value = 0.68877551020408101
$ws1.Range("B13").NumberFormat = "0.000000"
$ws1.Range("C13").NumberFormat = "0.000000"
$ws1.Range("F13").NumberFormat = "0.000000"

# --- Sheet 2: "Dep or Non-Dep With LOSOCV" ---
$ws2.Range("B6").Value = 0.57269064269064196
$ws2.Range("C6").Value = 0.32727272727272699
$ws2.Range("D6").Value = 0.10047619047619
$ws2.Range("E6").Value = 0.146961826052735
$ws2.Range("F6").Value = 0.57269064269064196
$ws2.Range("F6").NumberFormat = "0.000000"

$ws2.Range("B13").Value = 0.51526271708089799
$ws2.Range("C13").Value = 0.4
$ws2.Range("D13").Value = 0.135491932310114
$ws2.Range("E13").Value = 0.192604768968405
$ws2.Range("F13").Value = 0.51526271708089799
$ws2.Range("F13").NumberFormat = "0.000000"
